# Apply crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.357.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = "'1.859.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = "'314.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").Value = "'0.3715"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07313"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = "'0.8918"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("D11").Value = "'20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = "'0.07841"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.868.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'5.399"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").Value = "'6.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = "'91.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = "'0.000008946"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = "'14.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = "'27.393.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("D22").Value = "'5.135"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").Value = "'2.042.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").Value = "'1.925"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.10%  '
$ws.Range("D26").Value = "'152.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = "'18.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = "'2.058"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'116.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'5.093"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("D31").Value = "'0.08830"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = "'0.7729"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.79%  '
$ws.Range("D33").Value = "'3.075"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.65%  '
$ws.Range("D34").Value = "'1.176"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.75%  '
$ws.Range("D35").Value = "'4.517"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("D36").Value = "'2.732"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.75%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = "'1.080"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.01960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").Value = "'0.05269"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'2.971"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = "'7.037"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D42").Value = "'0.5145"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").Value = "'8.474"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.08%  '
$ws.Range("D45").Value = "'0.4804"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("D46").Value = "'10.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").Value = "'102.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = "'1.648"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.06220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").Value = "'65.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.46%  '
